# Daily attendance processing - 2025-11-06 23:42:54
#
# The "Recorded By" column (G) contains comma-separated lists of the
# users/systems that recorded a session. A couple of recurring value
# combinations need their entries reordered:
#
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com"     -> "System, backup@backdoor.com, system"
#
# Walk every used row on the active sheet and fix up column G wherever it
# exactly matches one of those two values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$colG = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, system, backup@backdoor.com") {
        $cell.Value = "System, backup@backdoor.com, system"
    }
}
